# Append: 2025-09-15 06:27 JST
# Updates the "ランサーズ" (Lancers) scrape sheet:
#   - refresh the "取得日時" (fetched-at) timestamp on every existing row
#   - roll rows 5-7 to the newest three scraped listings
#   - widen column B / narrow column D slightly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-15 06:27:55"

# --- Refresh "取得日時" for every data row (rows 2-7) ---
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- Row 5: new listing ---
$ws.Range("B5").Value = "【急募】Gasを使用した公式LINEチャットbotの作成依頼"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5393641"
$ws.Range("G5").Value = 118
$ws.Range("H5").Value = "★bot"

# --- Row 6: new listing ---
$ws.Range("B6").Value = "仮想通貨トレードの運用とコンサル【1名】のみ募集"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5393695"
$ws.Range("G6").Value = 55
$ws.Range("H6").Value = "◆コンサル"

# --- Row 7: new listing ---
$ws.Range("B7").Value = "【急募】トライアスロン大会運営支援システムの動作チェック、デバグ、品質確認業務委託費"
$ws.Range("D7").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5393606"
$ws.Range("G7").Value = 33

# --- Repoint the F5/F6/F7 hyperlinks at the new listing URLs ---
$ws.Hyperlinks.Item(4).Address = "https://www.lancers.jp/work/detail/5393641"
$ws.Hyperlinks.Item(5).Address = "https://www.lancers.jp/work/detail/5393695"
$ws.Hyperlinks.Item(6).Address = "https://www.lancers.jp/work/detail/5393606"

# --- Column width tweaks ---
$ws.Columns.Item(2).ColumnWidth = 43.17
$ws.Columns.Item(4).ColumnWidth = 27.17
